$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: change D2 from FAPs to ECs, update numeric values
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.359437
$ws.Range("H2").Value = 55.078311
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.09433999999999999
$ws.Range("N2").Value = 0.28302
$ws.Range("O2").Value = 0.05191071108246543
$ws.Range("P2").Value = 0.05191071108246543
$ws.Range("Q2").Value = 1.73202928658
$ws.Range("R2").Value = 15.58826357922
$ws.Range("S2").Value = 0.05191071108246543
$ws.Range("T2").Value = 0.05191071108246543

# Row 3: now represents FAPs target (previously sCs); update all values
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 18.359437
$ws.Range("H3").Value = 55.078311
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9431116666666667
$ws.Range("N3").Value = 2.829335
$ws.Range("O3").Value = 0.5189484550226392
$ws.Range("P3").Value = 0.5189484550226391
$ws.Range("Q3").Value = 17.31499922813167
$ws.Range("R3").Value = 155.834993053185
$ws.Range("S3").Value = 0.5189484550226392
$ws.Range("T3").Value = 0.5189484550226391

# Row 4: new row representing sCs target (moved from old row 3 with new values)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pf4"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 18.359437
$ws.Range("H4").Value = 55.078311
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7798996666666667
$ws.Range("N4").Value = 2.339699
$ws.Range("O4").Value = 0.4291408338948954
$ws.Range("P4").Value = 0.4291408338948954
$ws.Range("Q4").Value = 14.31851879648767
$ws.Range("R4").Value = 128.866669168389
$ws.Range("S4").Value = 0.4291408338948954
$ws.Range("T4").Value = 0.4291408338948954
